# Daily attendance processing - 2025-10-21 20:22:45
# Moves the exact-case token "System" to the front of the comma-separated
# "Recorded By" list (column G) wherever it appears in that list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CaseSensEq($a, $b) {
    if ($a.Length -ne $b.Length) { return $false }
    for ($i = 0; $i -lt $a.Length; $i++) {
        if ([int]$a[$i] -ne [int]$b[$i]) { return $false }
    }
    return $true
}

function Transform-RecordedBy($s) {
    if ([string]::IsNullOrEmpty($s)) { return $s }

    $parts = @($s -split ", ")

    $hasSystem = $false
    foreach ($p in $parts) {
        if (CaseSensEq $p "System") { $hasSystem = $true }
    }
    if (-not $hasSystem) { return $s }

    # Remove the first exact-case "System" entry, keep the rest in order.
    $rest = New-Object System.Collections.ArrayList
    $removed = $false
    foreach ($p in $parts) {
        if ((-not $removed) -and (CaseSensEq $p "System")) {
            $removed = $true
        } else {
            [void]$rest.Add($p)
        }
    }

    $new = @("System")
    foreach ($p in $rest) { $new += $p }

    return ($new -join ", ")
}

$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ([string]::IsNullOrEmpty($val)) { continue }

    $newVal = Transform-RecordedBy $val
    if ($newVal -cne $val) {
        $cell.Value = $newVal
    }
}
